$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D-column values are written as text (not auto-converted to numbers)
# by temporarily formatting the range as Text, then restoring the original style.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "21.411.47"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").Value = "1.552.58"
$ws.Range("E3").Value = "  +5.27%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "0.9696"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "283.37"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").Value = "0.3635"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "0.3223"
$ws.Range("E8").Value = "  +4.77%  "
$ws.Range("D9").Value = "41.11"
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("D11").Value = "0.06962"
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("D12").Value = "0.9987"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "5.739"
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").Value = "18.97"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("D15").Value = "6.433"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("D16").Value = "0.00001057"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "0.9675"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "1.544.29"
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("D19").Value = "0.06162"
$ws.Range("E19").Value = "  +4.79%  "
$ws.Range("D20").Value = "73.36"
$ws.Range("E20").Value = "  +5.76%  "
$ws.Range("D21").Value = "5.749"
$ws.Range("E21").Value = "  +5.24%  "
$ws.Range("D22").Value = "15.30"
$ws.Range("E22").Value = "  +5.94%  "
$ws.Range("D23").Value = "11.43"
$ws.Range("E23").Value = "  +4.19%  "
$ws.Range("D24").Value = "2.321"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").Value = "21.315.75"
$ws.Range("E25").Value = "  +3.99%  "
$ws.Range("D26").Value = "147.70"
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").Value = "2.282"
$ws.Range("E27").Value = "  +6.11%  "
$ws.Range("D28").Value = "17.91"
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").Value = "1.714.17"
$ws.Range("E29").Value = "  +5.27%  "
$ws.Range("D30").Value = "119.32"
$ws.Range("E30").Value = "  +4.79%  "
$ws.Range("D31").Value = "4.067"
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("D32").Value = "0.8755"
$ws.Range("E32").Value = "  +8.89%  "
$ws.Range("D33").Value = "5.287"
$ws.Range("E33").Value = "  +5.77%  "
$ws.Range("D34").Value = "0.08078"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").Value = "1.524"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "5.003"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").Value = "1.203"
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "0.05884"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").Value = "7.945"
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("D41").Value = "10.79"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "0.1934"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").Value = "0.9668"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "0.5513"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").Value = "12.66"
$ws.Range("E45").Value = "  +4.92%  "
$ws.Range("D46").Value = "3.576"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").Value = "0.5512"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").Value = "122.79"
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("D49").Value = "1.884"
$ws.Range("E49").Value = "  +6.05%  "
$ws.Range("D50").Value = "0.06629"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "70.21"
$ws.Range("E51").Value = "  +4.83%  "

# Restore the original (default) style on the D column so no stray number format lingers
$dRange.Style = "Normal"

